# Grading quiz 2, updated warm ups
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J: "Quiz 2" header + per-student scores (out of 35)
$ws.Range("J1").Value2 = "Quiz 2"

$ws.Range("J2").Formula  = "=37/35"
$ws.Range("J3").Formula  = "=34/35"
$ws.Range("J4").Formula  = "=38/35"
$ws.Range("J7").Formula  = "=22/35"
$ws.Range("J8").Formula  = "=35/35"
$ws.Range("J9").Formula  = "=38/35"
$ws.Range("J13").Formula = "=31/35"
$ws.Range("J14").Formula = "=34/35"
$ws.Range("J15").Formula = "=37/35"

# Move the active selection to J10, matching the updated view state
[void]$ws.Range("J10").Select()
